$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.009813969220977025
$ws.Cells.Item(2, 3).Value = 5.462182818544318
$ws.Cells.Item(2, 4).Value = 0.7261695939485373
$ws.Cells.Item(3, 2).Value = 1.192547057092041
$ws.Cells.Item(3, 3).Value = 3.281385277379687
$ws.Cells.Item(3, 4).Value = 0.2931506806686283
$ws.Cells.Item(4, 2).Value = 1.566835668388872
$ws.Cells.Item(4, 3).Value = 1.588326095274687
$ws.Cells.Item(4, 4).Value = 0.3960503260662407
$ws.Cells.Item(5, 2).Value = 2.071074682766095
$ws.Cells.Item(5, 3).Value = 4.698687301403841
$ws.Cells.Item(5, 4).Value = 0.5068839738218642
$ws.Cells.Item(6, 2).Value = 3.382064143243381
$ws.Cells.Item(6, 3).Value = 4.467802489710113
$ws.Cells.Item(6, 4).Value = 0.4411536902635889
$ws.Cells.Item(7, 2).Value = 4.270133868550261
$ws.Cells.Item(7, 3).Value = 3.568317967225127
$ws.Cells.Item(7, 4).Value = 0.5835960882406216
$ws.Cells.Item(8, 2).Value = 4.592517278619616
$ws.Cells.Item(8, 3).Value = 5.419508559878333
$ws.Cells.Item(8, 4).Value = 0.7461803148957802
$ws.Cells.Item(9, 2).Value = 6.070835942967911
$ws.Cells.Item(9, 3).Value = 2.416906371367066
$ws.Cells.Item(9, 4).Value = 0.4995478249762403
$ws.Cells.Item(10, 2).Value = 6.392804556578985
$ws.Cells.Item(10, 3).Value = 4.135841112573103
$ws.Cells.Item(10, 4).Value = 0.7231982363165215
$ws.Cells.Item(11, 2).Value = 6.470317699512129
$ws.Cells.Item(11, 3).Value = 4.557761353592925
$ws.Cells.Item(11, 4).Value = 0.647556435676761
$ws.Cells.Item(12, 2).Value = 10.7562053487111
$ws.Cells.Item(12, 3).Value = 4.479062987485866
$ws.Cells.Item(12, 4).Value = 0.4660855379397868
$ws.Cells.Item(13, 2).Value = 12.68008485475643
$ws.Cells.Item(13, 3).Value = 2.212618915762245
$ws.Cells.Item(13, 4).Value = 0.376145226932482
$ws.Cells.Item(14, 2).Value = 12.89195321882943
$ws.Cells.Item(14, 3).Value = 5.568986762060234
$ws.Cells.Item(14, 4).Value = 0.7945611276433515
$ws.Cells.Item(15, 2).Value = 13.06139946820941
$ws.Cells.Item(15, 3).Value = 4.956751005551578
$ws.Cells.Item(15, 4).Value = 0.2834593826952755
$ws.Cells.Item(16, 2).Value = 14.34838051144141
$ws.Cells.Item(16, 3).Value = 6.369142437324185
$ws.Cells.Item(16, 4).Value = 0.4960800443222256
$ws.Cells.Item(17, 2).Value = 14.42324679271286
$ws.Cells.Item(17, 3).Value = 4.717595153131771
$ws.Cells.Item(17, 4).Value = 0.4350432802088509
$ws.Cells.Item(18, 2).Value = 17.25847032585683
$ws.Cells.Item(18, 3).Value = 5.33113806225834
$ws.Cells.Item(18, 4).Value = 0.2734534515482324
$ws.Cells.Item(19, 2).Value = 19.92206171700013
$ws.Cells.Item(19, 3).Value = 1.71514571030078
$ws.Cells.Item(19, 4).Value = 0.2926036960548201
$ws.Cells.Item(20, 2).Value = 23.08278072467364
$ws.Cells.Item(20, 3).Value = 5.14724165512163
$ws.Cells.Item(20, 4).Value = 0.2670374852960377
$ws.Cells.Item(21, 2).Value = 24.18989236034563
$ws.Cells.Item(21, 3).Value = 1.273163660290217
$ws.Cells.Item(21, 4).Value = 0.2375756593480204
$ws.Cells.Item(22, 2).Value = 24.48933441279039
$ws.Cells.Item(22, 3).Value = 5.638049567461533
$ws.Cells.Item(22, 4).Value = 1.007690996895927
$ws.Cells.Item(23, 2).Value = 25.92202583030161
$ws.Cells.Item(23, 3).Value = 5.017202498399529
$ws.Cells.Item(23, 4).Value = 0.4600187139217956
$ws.Cells.Item(24, 2).Value = 27.61918058073666
$ws.Cells.Item(24, 3).Value = 7.128937344834235
$ws.Cells.Item(24, 4).Value = 0.6305989415861005
$ws.Cells.Item(25, 2).Value = 27.6840127613614
$ws.Cells.Item(25, 3).Value = 3.555288353141759
$ws.Cells.Item(25, 4).Value = 0.4081254402198048
$ws.Cells.Item(26, 2).Value = 28.64010623458799
$ws.Cells.Item(26, 3).Value = 4.367303771870753
$ws.Cells.Item(26, 4).Value = 0.2201182784334469
$ws.Cells.Item(27, 2).Value = 29.27700640614891
$ws.Cells.Item(27, 3).Value = 2.561774208989367
$ws.Cells.Item(27, 4).Value = 0.2845284103600966
$ws.Cells.Item(28, 2).Value = 29.71630388974771
$ws.Cells.Item(28, 3).Value = 4.463789292588459
$ws.Cells.Item(28, 4).Value = 0.6464083941813835
$ws.Cells.Item(29, 2).Value = 30.0475402957403
$ws.Cells.Item(29, 3).Value = 5.049260323065087
$ws.Cells.Item(29, 4).Value = 0.6324599749413289
$ws.Cells.Item(30, 2).Value = 31.752263083003
$ws.Cells.Item(30, 3).Value = 6.798164295952496
$ws.Cells.Item(30, 4).Value = 0.5951790032988106
$ws.Cells.Item(31, 2).Value = 32.822750275225
$ws.Cells.Item(31, 3).Value = 3.145257121387247
$ws.Cells.Item(31, 4).Value = 0.3932055669661726
$ws.Cells.Item(32, 2).Value = 35.26427629289751
$ws.Cells.Item(32, 3).Value = 7.316554550625701
$ws.Cells.Item(32, 4).Value = 0.77679041732841
$ws.Cells.Item(33, 2).Value = 36.21812152556571
$ws.Cells.Item(33, 3).Value = 6.206922892595349
$ws.Cells.Item(33, 4).Value = 0.6002808179384113
$ws.Cells.Item(34, 2).Value = 36.41119441319319
$ws.Cells.Item(34, 3).Value = 10.20184830752223
$ws.Cells.Item(34, 4).Value = 0.4732458405060798
$ws.Cells.Item(35, 2).Value = 40.23090346143096
$ws.Cells.Item(35, 3).Value = 7.115512409165457
$ws.Cells.Item(35, 4).Value = 0.6417775522272048
$ws.Cells.Item(36, 2).Value = 40.6031106991014
$ws.Cells.Item(36, 3).Value = 9.941265617257175
$ws.Cells.Item(36, 4).Value = 0.3768509962478094
$ws.Cells.Item(37, 2).Value = 41.06000397794802
$ws.Cells.Item(37, 3).Value = 6.072939965884224
$ws.Cells.Item(37, 4).Value = 0.6184282011598021
$ws.Cells.Item(38, 2).Value = 41.77329725711184
$ws.Cells.Item(38, 3).Value = 5.003389327004911
$ws.Cells.Item(38, 4).Value = 0.3864879271982016
$ws.Cells.Item(39, 2).Value = 43.27566297447875
$ws.Cells.Item(39, 3).Value = 4.477937950478132
$ws.Cells.Item(39, 4).Value = 0.6930632522428677
$ws.Cells.Item(40, 2).Value = 43.52687972704184
$ws.Cells.Item(40, 3).Value = 9.224090278083894
$ws.Cells.Item(40, 4).Value = 0.5445991040934469
$ws.Cells.Item(41, 2).Value = 44.44644819957989
$ws.Cells.Item(41, 3).Value = 1.753745344624981
$ws.Cells.Item(41, 4).Value = 0.2996719965052247
$ws.Cells.Item(42, 2).Value = 44.8961634975113
$ws.Cells.Item(42, 3).Value = 5.935448748226894
$ws.Cells.Item(42, 4).Value = 0.2715471464536222
$ws.Cells.Item(43, 2).Value = 47.73511288424801
$ws.Cells.Item(43, 3).Value = 1.191993140685253
$ws.Cells.Item(43, 4).Value = 0.2321585447971928
$ws.Cells.Item(44, 2).Value = 49.92633862899439
$ws.Cells.Item(44, 3).Value = 8.242937729577966
$ws.Cells.Item(44, 4).Value = 0.6731898588032115
$ws.Cells.Item(45, 2).Value = 50.82201394110999
$ws.Cells.Item(45, 3).Value = 3.58463549528678
$ws.Cells.Item(45, 4).Value = 0.3449359193436335
$ws.Cells.Item(46, 2).Value = 52.86330516856284
$ws.Cells.Item(46, 3).Value = 3.01604934851022
$ws.Cells.Item(46, 4).Value = 0.29465670260001
$ws.Cells.Item(47, 2).Value = 53.24158689529139
$ws.Cells.Item(47, 3).Value = 10.52710352419931
$ws.Cells.Item(47, 4).Value = 0.624120814743659
$ws.Cells.Item(48, 2).Value = 54.97810273519707
$ws.Cells.Item(48, 3).Value = 8.472964529552129
$ws.Cells.Item(48, 4).Value = 0.4236995365752461
$ws.Cells.Item(49, 2).Value = 55.72446547011131
$ws.Cells.Item(49, 3).Value = 4.003522222932123
$ws.Cells.Item(49, 4).Value = 0.2756575128965081
$ws.Cells.Item(50, 2).Value = 56.17770280497886
$ws.Cells.Item(50, 3).Value = 6.837436478097663
$ws.Cells.Item(50, 4).Value = 0.477820256068485
$ws.Cells.Item(51, 2).Value = 56.57581747646673
$ws.Cells.Item(51, 3).Value = 6.240471964965037
$ws.Cells.Item(51, 4).Value = 0.7211646398844105
$ws.Cells.Item(52, 2).Value = 60.41808651136485
$ws.Cells.Item(52, 3).Value = 3.659982506470986
$ws.Cells.Item(52, 4).Value = 0.2429821834532013
$ws.Cells.Item(53, 2).Value = 64.98665520275823
$ws.Cells.Item(53, 3).Value = 3.721796058797615
$ws.Cells.Item(53, 4).Value = 0.4457190948507765
$ws.Cells.Item(54, 2).Value = 65.10725774047822
$ws.Cells.Item(54, 3).Value = 2.149021594718868
$ws.Cells.Item(54, 4).Value = 0.4561871334304363
$ws.Cells.Item(55, 2).Value = 65.25237543653881
$ws.Cells.Item(55, 3).Value = 8.145955335394975
$ws.Cells.Item(55, 4).Value = 0.5694979589469783
$ws.Cells.Item(56, 2).Value = 65.98981100699106
$ws.Cells.Item(56, 3).Value = 3.856927269729924
$ws.Cells.Item(56, 4).Value = 0.2922003894904837
$ws.Cells.Item(57, 2).Value = 69.98679088817806
$ws.Cells.Item(57, 3).Value = 3.118918085703063
$ws.Cells.Item(57, 4).Value = 0.5977389253796682
$ws.Cells.Item(58, 2).Value = 71.31492148984466
$ws.Cells.Item(58, 3).Value = 2.233744311033993
$ws.Cells.Item(58, 4).Value = 0.2229845595858626
$ws.Cells.Item(59, 2).Value = 72.53831015740562
$ws.Cells.Item(59, 3).Value = 8.903793297634307
$ws.Cells.Item(59, 4).Value = 0.4084051103098985
$ws.Cells.Item(60, 2).Value = 74.40993033280505
$ws.Cells.Item(60, 3).Value = 3.654097356018678
$ws.Cells.Item(60, 4).Value = 0.6242648973857674
$ws.Cells.Item(61, 2).Value = 74.75192375760891
$ws.Cells.Item(61, 3).Value = 1.382320021955507
$ws.Cells.Item(61, 4).Value = 0.3560040341798285
$ws.Cells.Item(62, 2).Value = 74.93909880401526
$ws.Cells.Item(62, 3).Value = 5.540263116033993
$ws.Cells.Item(62, 4).Value = 0.6019402940453928
$ws.Cells.Item(63, 2).Value = 75.48269369780367
$ws.Cells.Item(63, 3).Value = 4.673724270885131
$ws.Cells.Item(63, 4).Value = 0.4685739852307593
$ws.Cells.Item(64, 2).Value = 77.75151709634163
$ws.Cells.Item(64, 3).Value = 4.542762091253021
$ws.Cells.Item(64, 4).Value = 0.4461267402401625
$ws.Cells.Item(65, 2).Value = 80.13149947336636
$ws.Cells.Item(65, 3).Value = 2.60875082246766
$ws.Cells.Item(65, 4).Value = 0.2937017422531613
$ws.Cells.Item(66, 2).Value = 81.183438795114
$ws.Cells.Item(66, 3).Value = 5.432632010864847
$ws.Cells.Item(66, 4).Value = 0.4075256259236966
$ws.Cells.Item(67, 2).Value = 81.27616644013899
$ws.Cells.Item(67, 3).Value = 3.604367172939565
$ws.Cells.Item(67, 4).Value = 0.4670275678677585
$ws.Cells.Item(68, 2).Value = 83.15376327290149
$ws.Cells.Item(68, 3).Value = 4.107960174013623
$ws.Cells.Item(68, 4).Value = 0.2066852089862372
$ws.Cells.Item(69, 2).Value = 83.39713331213849
$ws.Cells.Item(69, 3).Value = 4.006240592700196
$ws.Cells.Item(69, 4).Value = 0.747149462053189
$ws.Cells.Item(70, 2).Value = 85.47273514521804
$ws.Cells.Item(70, 3).Value = 5.629065470825552
$ws.Cells.Item(70, 4).Value = 0.6301387307983135
$ws.Cells.Item(71, 2).Value = 86.76304712078324
$ws.Cells.Item(71, 3).Value = 2.892404332343843
$ws.Cells.Item(71, 4).Value = 0.3883951326189075
$ws.Cells.Item(72, 2).Value = 87.40772906464281
$ws.Cells.Item(72, 3).Value = 10.46778992347473
$ws.Cells.Item(72, 4).Value = 0.4858406401008161
$ws.Cells.Item(73, 2).Value = 92.93616654884543
$ws.Cells.Item(73, 3).Value = 7.02454441327987
$ws.Cells.Item(73, 4).Value = 0.5481133999659223
$ws.Cells.Item(74, 2).Value = 96.06539287182014
$ws.Cells.Item(74, 3).Value = 2.05179919860438
$ws.Cells.Item(74, 4).Value = 0.3053026253373703

$ws.Range("A75:D81").EntireRow.Delete()
